$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- helper: write a date-looking string as TEXT, not as an auto-converted
# Excel date serial. We stage it in a scratch cell that has been forced to
# Text format, then copy *values only* onto the real target so the target
# cell's own style/number-format stays untouched (matches the source file,
# where these are plain shared-string cells with no special numFmt).
$scratch = $ws.Range("Z100")
$scratch.NumberFormat = "@"
$scratch.Value = "10-11-2024"
$scratch.Copy()

function Set-TextDate($cell) {
    $cell.PasteSpecial(-4163)   # xlPasteValues
}

# --- Row 2-5: method names changed (new "contact info / registration" flow)
# and the execution date moves from 04-11-2024 to 10-11-2024. "Passed" (col B)
# and its green highlight fill are unchanged.
$ws.Range("A2").Value = "verifyCustomerNavigationToRegistrationPage"
Set-TextDate($ws.Range("C2"))

$ws.Range("A3").Value = "verifyNewCustomerRegistrationSubmissionFlow"
Set-TextDate($ws.Range("C3"))

$ws.Range("A4").Value = "verifyCustomerRegistrationAndLoginNavigation"
Set-TextDate($ws.Range("C4"))

$ws.Range("A5").Value = "verifyCustomerEmailActivation"
Set-TextDate($ws.Range("C5"))

# --- New rows 6-11: continuation of the existing suite plus the new flow.
$newRows = @(
    "verifyCustomerSuccessfulLogin",
    "verifyCustomerNavigationAfterLogin",
    "verifyCustomerPreferredPackageSelection",
    "verifyCustomerNavigationAfterSaving",
    "verifyCustomerBasicInfoEntry",
    "verifyCustomerNavigationAfterSaving"
)

$rowIndex = 6
foreach ($methodName in $newRows) {
    $ws.Cells.Item($rowIndex, 1).Value = $methodName

    # Col B: "Passed" with the same green fill used by the existing rows -
    # copy format+value from the row-2 template cell so it matches exactly.
    $ws.Range("B2").Copy()
    $targetCell = $ws.Cells.Item($rowIndex, 2)
    $targetCell.PasteSpecial(-4104)  # xlPasteAll

    Set-TextDate($ws.Cells.Item($rowIndex, 3))

    $rowIndex = $rowIndex + 1
}

$scratch.Clear()

# --- Column A widens to fit the longer method names now in use.
$ws.Columns.Item(1).ColumnWidth = 43.1

Write-Host "done"
